# Update the "想去人数" (want-to-go count) figures that changed between
# successive generations of the scraped data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1582
$ws1.Range("F7").Value = 6447
$ws1.Range("F11").Value = 5744
$ws1.Range("F17").Value = 77
$ws1.Range("F21").Value = 327
$ws1.Range("F24").Value = 4094
$ws1.Range("F26").Value = 180

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1582
$ws4.Range("F7").Value = 6447
$ws4.Range("F11").Value = 5744
$ws4.Range("F17").Value = 77
$ws4.Range("F21").Value = 327
$ws4.Range("F24").Value = 4094
$ws4.Range("F27").Value = 180
